$wb = $excel.ActiveWorkbook

# --- Sheet "10.10 - 16.10." (index 4): fill previously empty row 6, and set explicit 0s in row 8 ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("C6").Value = 2
$ws4.Range("D6").Value = 1
$ws4.Range("E6").Value = 0
$ws4.Range("F6").Value = 1
$ws4.Range("G6").Value = 2
$ws4.Range("H6").Value = 0
$ws4.Range("I6").Value = 1
$ws4.Range("E8").Value = 0
$ws4.Range("H8").Value = 0

# --- Sheet "17.10. - 23.10." (index 5): fill/modify rows 4-8 ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("C4").Value = 3
$ws5.Range("D4").Value = 1
$ws5.Range("G4").Value = 3

$ws5.Range("D5").Value = 3
$ws5.Range("E5").Value = 0
$ws5.Range("F5").Value = 0
$ws5.Range("G5").Value = 4
$ws5.Range("H5").Value = 4
$ws5.Range("I5").Value = 5

$ws5.Range("C6").Value = 2
$ws5.Range("D6").Value = 1

$ws5.Range("C7").Value = 1
$ws5.Range("D7").Value = 2
$ws5.Range("E7").Value = 0
$ws5.Range("F7").Value = 5
$ws5.Range("G7").Value = 3
$ws5.Range("H7").Value = 1
$ws5.Range("I7").Value = 3

$ws5.Range("C8").Value = 3
$ws5.Range("D8").Value = 2
$ws5.Range("E8").Value = 0
$ws5.Range("F8").Value = 3
$ws5.Range("G8").Value = 3
$ws5.Range("H8").Value = 2
$ws5.Range("I8").Value = 3

# --- View / selection state ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Activate()
$ws3.Range("C8").Select()

$ws4.Activate()
$ws4.Range("F13").Select()

$ws5.Activate()
$ws5.Range("C5:I5").Select()
